# Apply crypto price/volume updates per the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.573.23"
$ws.Range("E2").Value = "'  +2.48%  "
$ws.Range("D3").Value = "'1.673.18"
$ws.Range("E3").Value = "'  +2.15%  "
$ws.Range("D4").Value = "'0.9988"
$ws.Range("E4").Value = "'  +0.24%  "
$ws.Range("D5").Value = "'239.85"
$ws.Range("E5").Value = "'  +1.32%  "
$ws.Range("E6").Value = "'  +0.00%  "
$ws.Range("D7").Value = "'0.4771"
$ws.Range("E7").Value = "'  +0.96%  "
$ws.Range("D8").Value = "'0.2628"
$ws.Range("E8").Value = "'  +2.54%  "
$ws.Range("D9").Value = "'0.06181"
$ws.Range("E9").Value = "'  +2.78%  "
$ws.Range("D10").Value = "'1.673.26"
$ws.Range("E10").Value = "'  +2.14%  "
$ws.Range("D11").Value = "'0.06994"
$ws.Range("E11").Value = "'  -1.05%  "
$ws.Range("D12").Value = "'14.90"
$ws.Range("E12").Value = "'  +0.47%  "
$ws.Range("D13").Value = "'0.5925"
$ws.Range("E13").Value = "'  -3.92%  "
$ws.Range("D14").Value = "'4.387"
$ws.Range("E14").Value = "'  +0.33%  "
$ws.Range("D15").Value = "'75.46"
$ws.Range("E15").Value = "'  +3.66%  "
$ws.Range("D16").Value = "'0.9997"
$ws.Range("D17").Value = "'0.9993"
$ws.Range("E17").Value = "'  +0.17%  "
$ws.Range("D18").Value = "'25.564.52"
$ws.Range("E18").Value = "'  +2.48%  "
$ws.Range("D19").Value = "'0.000006765"
$ws.Range("E19").Value = "'  +2.76%  "
$ws.Range("D20").Value = "'11.45"
$ws.Range("E20").Value = "'  +2.45%  "
$ws.Range("D21").Value = "'1.887.26"
$ws.Range("E21").Value = "'  +2.30%  "
$ws.Range("D22").Value = "'4.460"
$ws.Range("E22").Value = "'  +1.02%  "
$ws.Range("D23").Value = "'8.787"
$ws.Range("E23").Value = "'  +2.14%  "
$ws.Range("D24").Value = "'5.283"
$ws.Range("E24").Value = "'  +0.05%  "
$ws.Range("D25").Value = "'136.97"
$ws.Range("E25").Value = "'  +2.94%  "
$ws.Range("E26").Value = "'  +1.70%  "
$ws.Range("D27").Value = "'1.387"
$ws.Range("E27").Value = "'  +1.88%  "
$ws.Range("D28").Value = "'1.732"
$ws.Range("E28").Value = "'  +4.34%  "
$ws.Range("D29").Value = "'104.79"
$ws.Range("E29").Value = "'  +2.12%  "
$ws.Range("D30").Value = "'3.974"
$ws.Range("E30").Value = "'  +5.92%  "
$ws.Range("D31").Value = "'0.07865"
$ws.Range("E31").Value = "'  +1.73%  "
$ws.Range("D32").Value = "'3.645"
$ws.Range("E32").Value = "'  +2.45%  "
$ws.Range("E33").Value = "'  +0.03%  "
$ws.Range("D34").Value = "'0.04285"
$ws.Range("E34").Value = "'  -0.99%  "
$ws.Range("E35").Value = "'  +0.87%  "
$ws.Range("D36").Value = "'0.9591"
$ws.Range("E36").Value = "'  +3.97%  "
$ws.Range("D37").Value = "'0.6084"
$ws.Range("E37").Value = "'  +4.47%  "
$ws.Range("D38").Value = "'2.598"
$ws.Range("E38").Value = "'  +1.06%  "
$ws.Range("D39").Value = "'0.8930"
$ws.Range("E39").Value = "'  +8.30%  "
$ws.Range("D40").Value = "'0.9996"
$ws.Range("D41").Value = "'1.868"
$ws.Range("E41").Value = "'  +3.73%  "
$ws.Range("D42").Value = "'0.01487"
$ws.Range("E42").Value = "'  -4.41%  "
$ws.Range("E43").Value = "'  -1.20%  "
$ws.Range("D44").Value = "'0.3768"
$ws.Range("E44").Value = "'  +1.33%  "
$ws.Range("D45").Value = "'4.915"
$ws.Range("E45").Value = "'  +3.75%  "
$ws.Range("D46").Value = "'0.1122"
$ws.Range("E46").Value = "'  +1.36%  "
$ws.Range("D47").Value = "'6.238"
$ws.Range("E47").Value = "'  +2.42%  "
$ws.Range("D48").Value = "'0.05268"
$ws.Range("E48").Value = "'  +0.96%  "
$ws.Range("E49").Value = "'  +1.23%  "
$ws.Range("D50").Value = "'7.463"
$ws.Range("E50").Value = "'  +4.24%  "
$ws.Range("B51").Value = "'NEARProtocol"
$ws.Range("C51").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.210"
$ws.Range("E51").Value = "'  +2.65%  "
